$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: "..._old" -> "..._FV2410", "..._new" -> "..._FV2504" ---
# (the lone "diff" header in between has neither suffix and is left untouched)
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value2
    if ($header -match '_old$') {
        $cell.Value = ($header -replace '_old$', '_FV2410')
    } elseif ($header -match '_new$') {
        $cell.Value = ($header -replace '_new$', '_FV2504')
    }
}

# --- Turn the data range into an Excel Table (ListObject) ---
$range = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- Freeze the header row (split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
